$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (C, D, H) ---
# Excel ColumnWidth has a fixed +0.8333333333333334 offset vs. the stored OOXML width
$ws.Columns.Item(3).ColumnWidth = 65.16666666666667   # -> stored width 66
$ws.Columns.Item(4).ColumnWidth = 56.166666666666664  # -> stored width 57
$ws.Columns.Item(8).ColumnWidth = 32.166666666666664  # -> stored width 33

# --- Highlight PREMIUM = Yes cells in column E with yellow fill ---
# Do this first so the new fill/style is created as cellXfs index 3,
# matching the order in which the workbook introduces it.
# Build a single multi-area range so only one new fill/style is created
$premiumYes = $excel.Union($ws.Range("E2"), $ws.Range("E3"), $ws.Range("E5"), $ws.Range("E6"), $ws.Range("E15"), $ws.Range("E16"))
foreach ($area in $premiumYes.Areas) {
    $area.Interior.Color = 65535   # RGB yellow (255,255,0)
}

# --- Write data rows 2-16 ---
# Row 2
$idCell = $ws.Cells.Item(2, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1297124'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(2, 2).Value = 'https://aiesec.org/opportunity/global-talent/1297124'
$ws.Cells.Item(2, 3).Value = 'ACE Program | Portuguese Talent Acquisition Specialist'
$ws.Cells.Item(2, 4).Value = 'Chennai, Tamil Nadu, India'
$ws.Cells.Item(2, 5).Value = 'Yes'
$ws.Cells.Item(2, 6).Value = '12 applicants'
$ws.Cells.Item(2, 7).Value = '6 - 18 Months'
$ws.Cells.Item(2, 8).Value = 'Tata Consultancy Services Ltd.'

# Row 3
$idCell = $ws.Cells.Item(3, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1327965'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(3, 2).Value = 'https://aiesec.org/opportunity/global-talent/1327965'
$ws.Cells.Item(3, 3).Value = 'ACE Program | German Financial Analyst'
$ws.Cells.Item(3, 4).Value = 'Thane, Maharashtra, India'
$ws.Cells.Item(3, 5).Value = 'Yes'
$ws.Cells.Item(3, 6).Value = '12 applicants'
$ws.Cells.Item(3, 7).Value = '6 - 18 Months'
$ws.Cells.Item(3, 8).Value = 'Tata Consultancy Services Ltd.'

# Row 4
$idCell = $ws.Cells.Item(4, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1330684'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(4, 2).Value = 'https://aiesec.org/opportunity/global-talent/1330684'
$ws.Cells.Item(4, 3).Value = 'Account Management Business Process Analyst'
$ws.Cells.Item(4, 4).Value = 'Panamá, Provincia de Panamá, Panamá'
$ws.Cells.Item(4, 5).Value = 'No'
$ws.Cells.Item(4, 6).Value = '2 applicants'
$ws.Cells.Item(4, 7).Value = '6 - 18 Months'
$ws.Cells.Item(4, 8).Value = 'G4S Panamá'

# Row 5
$idCell = $ws.Cells.Item(5, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1330679'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(5, 2).Value = 'https://aiesec.org/opportunity/global-talent/1330679'
$ws.Cells.Item(5, 3).Value = 'ACE Program | Global HR & RMG Business Partner (AIESECers Only)'
$ws.Cells.Item(5, 4).Value = 'Hyderabad, Telangana, India'
$ws.Cells.Item(5, 5).Value = 'Yes'
$ws.Cells.Item(5, 6).Value = '1 applicant'
$ws.Cells.Item(5, 7).Value = '6 - 18 Months'
$ws.Cells.Item(5, 8).Value = 'Tata Consultancy Services Ltd.'

# Row 6
$idCell = $ws.Cells.Item(6, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1330678'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(6, 2).Value = 'https://aiesec.org/opportunity/global-talent/1330678'
$ws.Cells.Item(6, 3).Value = 'ACE Program | Global Coordinator (AIESECers Only)'
$ws.Cells.Item(6, 4).Value = 'Hyderabad, Telangana, India'
$ws.Cells.Item(6, 5).Value = 'Yes'
$ws.Cells.Item(6, 6).Value = '1 applicant'
$ws.Cells.Item(6, 7).Value = '6 - 18 Months'
$ws.Cells.Item(6, 8).Value = 'Tata Consultancy Services Ltd.'

# Row 7
$idCell = $ws.Cells.Item(7, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1330676'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(7, 2).Value = 'https://aiesec.org/opportunity/global-talent/1330676'
$ws.Cells.Item(7, 3).Value = 'Marketing Analyst'
$ws.Cells.Item(7, 4).Value = 'Panamá, Provincia de Panamá, Panamá'
$ws.Cells.Item(7, 5).Value = 'No'
$ws.Cells.Item(7, 6).Value = '2 applicants'
$ws.Cells.Item(7, 7).Value = '6 - 18 Months'
$ws.Cells.Item(7, 8).Value = 'GSK - GlaxoSmithKline'

# Row 8
$idCell = $ws.Cells.Item(8, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1330246'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(8, 2).Value = 'https://aiesec.org/opportunity/global-talent/1330246'
$ws.Cells.Item(8, 3).Value = 'Mobile Developer'
$ws.Cells.Item(8, 4).Value = 'Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt'
$ws.Cells.Item(8, 5).Value = 'No'
$ws.Cells.Item(8, 6).Value = '3 applicants'
$ws.Cells.Item(8, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(8, 8).Value = 'Safarni'

# Row 9
$idCell = $ws.Cells.Item(9, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1328367'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(9, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328367'
$ws.Cells.Item(9, 3).Value = 'Interior Designer'
$ws.Cells.Item(9, 4).Value = 'Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt'
$ws.Cells.Item(9, 5).Value = 'No'
$ws.Cells.Item(9, 6).Value = '5 applicants'
$ws.Cells.Item(9, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(9, 8).Value = 'ASG Engineering'

# Row 10
$idCell = $ws.Cells.Item(10, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1328365'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(10, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328365'
$ws.Cells.Item(10, 3).Value = 'Graphic Designer'
$ws.Cells.Item(10, 4).Value = 'Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt'
$ws.Cells.Item(10, 5).Value = 'No'
$ws.Cells.Item(10, 6).Value = '7 applicants'
$ws.Cells.Item(10, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(10, 8).Value = 'ASG Engineering'

# Row 11
$idCell = $ws.Cells.Item(11, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1328363'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(11, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328363'
$ws.Cells.Item(11, 3).Value = 'Marketing Specialist'
$ws.Cells.Item(11, 4).Value = 'Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt'
$ws.Cells.Item(11, 5).Value = 'No'
$ws.Cells.Item(11, 6).Value = '15 applicants'
$ws.Cells.Item(11, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(11, 8).Value = 'ASG Engineering'

# Row 12
$idCell = $ws.Cells.Item(12, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1328345'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(12, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328345'
$ws.Cells.Item(12, 3).Value = 'Design Engineer'
$ws.Cells.Item(12, 4).Value = 'Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt'
$ws.Cells.Item(12, 5).Value = 'No'
$ws.Cells.Item(12, 6).Value = '9 applicants'
$ws.Cells.Item(12, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(12, 8).Value = 'print shop'

# Row 13
$idCell = $ws.Cells.Item(13, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1327904'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(13, 2).Value = 'https://aiesec.org/opportunity/global-talent/1327904'
$ws.Cells.Item(13, 3).Value = 'UX Research Trainee'
$ws.Cells.Item(13, 4).Value = 'Bruxelles, Belgio'
$ws.Cells.Item(13, 5).Value = 'No'
$ws.Cells.Item(13, 6).Value = '132 applicants'
$ws.Cells.Item(13, 7).Value = '6 - 18 Months'
$ws.Cells.Item(13, 8).Value = 'UCB'

# Row 14
$idCell = $ws.Cells.Item(14, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1326473'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(14, 2).Value = 'https://aiesec.org/opportunity/global-talent/1326473'
$ws.Cells.Item(14, 3).Value = 'Business Development Intern'
$ws.Cells.Item(14, 4).Value = 'Pune, India'
$ws.Cells.Item(14, 5).Value = 'No'
$ws.Cells.Item(14, 6).Value = '0 applicants'
$ws.Cells.Item(14, 7).Value = '3 - 6 Months'
$ws.Cells.Item(14, 8).Value = 'Spidron Tech LLP'

# Row 15
$idCell = $ws.Cells.Item(15, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1326041'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(15, 2).Value = 'https://aiesec.org/opportunity/global-talent/1326041'
$ws.Cells.Item(15, 3).Value = 'ACE Program | Spanish Financial Analyst'
$ws.Cells.Item(15, 4).Value = 'Thane, Maharashtra, India'
$ws.Cells.Item(15, 5).Value = 'Yes'
$ws.Cells.Item(15, 6).Value = '36 applicants'
$ws.Cells.Item(15, 7).Value = '6 - 18 Months'
$ws.Cells.Item(15, 8).Value = 'Tata Consultancy Services Ltd.'

# Row 16
$idCell = $ws.Cells.Item(16, 1)
$idCell.NumberFormat = "@"        # force text so the numeric-looking ID is not converted to a number
$idCell.Value = '1305153'
$idCell.Style = "Normal"          # drop the temporary text-format style again, keep the stored text value
$ws.Cells.Item(16, 2).Value = 'https://aiesec.org/opportunity/global-talent/1305153'
$ws.Cells.Item(16, 3).Value = 'ACE Program | Spanish Talent Acquisition Specialist'
$ws.Cells.Item(16, 4).Value = 'Chennai, Tamil Nadu, India'
$ws.Cells.Item(16, 5).Value = 'Yes'
$ws.Cells.Item(16, 6).Value = '56 applicants'
$ws.Cells.Item(16, 7).Value = '6 - 18 Months'
$ws.Cells.Item(16, 8).Value = 'Tata Consultancy Services Ltd.'

Write-Host "Edit applied"
